# Apply "Alteração para diferentes tipos de bimestrais" edits
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "-"
$ws.Range("D2").Value = "['MCT-3A-Tecnologia da soldagem', 'MCT-3A-Tecnologia da soldagem', -, -]"

$ws.Range("B3").Value = "['MEC-3B-Tec. Soldagem', -, 'MEC-3B-Tec. Soldagem', -]"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = "[-, 'MCT-3A-Tecnologia da soldagem', 'MCT-3A-Tecnologia da soldagem', -]"

$ws.Range("B4").Value = "['MEC-3B-Tec. Soldagem', -, 'MEC-3B-Tec. Soldagem', -]"
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = "-"

$ws.Range("B6").Value = "-"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = "-"

$ws.Range("D7").Value = "-"

$ws.Range("C19").Value = "['MEC-2NA-Soldagem', 'MEC-2NA-Soldagem', 'MEC-2NA-Soldagem', 'MEC-2NA-Soldagem']"
$ws.Range("F19").Value = "ELM-1NA-Gest. Int."

$ws.Range("B20").Value = "MEC-2NA-Gest. Int."
$ws.Range("E20").Value = "MEC-2NB-Gest. Int."
$ws.Range("F20").Value = "ELM-1NA-Gest. Int."

$ws.Range("B21").Value = "MEC-2NA-Gest. Int."
$ws.Range("E21").Value = "MEC-2NB-Gest. Int."

$wb.Save()
